$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting existing rows 10-26 down to 11-27.
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with the "Sapphire" record.
$ws.Cells.Item(10, 1).Value = 1
$ws.Cells.Item(10, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(10, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(10, 4).Value = 44966
$ws.Cells.Item(10, 5).Value = 15
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100103
$ws.Cells.Item(10, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(10, 9).Value = 100103002
$ws.Cells.Item(10, 10).Value = "Ciruela"
$ws.Cells.Item(10, 11).Value = "Sapphire"
$ws.Cells.Item(10, 12).Value = "Segunda"
$ws.Cells.Item(10, 13).Value = 250
$ws.Cells.Item(10, 14).Value = 18000
$ws.Cells.Item(10, 15).Value = 20000
$ws.Cells.Item(10, 16).Value = 19000
$ws.Cells.Item(10, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(10, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(10, 19).Value = 1056
$ws.Cells.Item(10, 20).Value = 18

# Apply the same date-time display format as the other "Fecha" column cells
# (the inserted row inherits this from row 9 via the Insert shift already,
# but set it explicitly to be safe).
$ws.Cells.Item(10, 4).NumberFormat = $ws.Cells.Item(11, 4).NumberFormat
